# Apply newest airtoxics NATA data updates to the allocation rule summary tables.
$wb = $excel.ActiveWorkbook

# --- "Means" sheet ---
$wsMeans = $wb.Worksheets.Item("Means")

# Row 9: Total Cancer Risk (per million)
$wsMeans.Range("B9").Value = 26
$wsMeans.Range("C9").Value = 28
$wsMeans.Range("D9").Value = 17
$wsMeans.Range("E9").Value = 19
$wsMeans.Range("F9").Value = 19
$wsMeans.Range("G9").Value = 19

# Row 10: Total Respiratory (hazard quotient)
$wsMeans.Range("B10").Value = 0.31
$wsMeans.Range("C10").Value = 0.3
$wsMeans.Range("F10").Value = 0.18
$wsMeans.Range("G10").Value = 0.18

# --- "Standard Deviations" sheet ---
$wsSD = $wb.Worksheets.Item("Standard Deviations")

# Row 9: Total Cancer Risk (per million) SD
$wsSD.Range("B9").Value = 8.3
$wsSD.Range("C9").Value = 13
$wsSD.Range("D9").Value = 5.8
$wsSD.Range("E9").Value = 3.2
$wsSD.Range("F9").Value = 3
$wsSD.Range("G9").Value = 3.3

# Row 10: Total Respiratory (hazard quotient) SD
$wsSD.Range("B10").Value = 0.11
$wsSD.Range("C10").Value = 0.11
$wsSD.Range("F10").Value = 0.036
$wsSD.Range("G10").Value = 0.039

$wb.Save()
